# edit.ps1
# Applies the changes described by the diff:
#  - Renames several header labels (shared strings) in row 1
#  - Updates the "Gross_National_Income"/GDP (column C) numeric values for all data rows
#  - Flips three "Colony" (column AL) flags from 0 to 1 for the Haiti rows (12, 28, 40)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label renames (row 1) ---
$ws.Range("C1").Value  = "GDP"
$ws.Range("E1").Value  = "Budget_Previous_Year"
$ws.Range("F1").Value  = "LatinAmerica"
$ws.Range("G1").Value  = "Africa"
$ws.Range("H1").Value  = "Confessional"
$ws.Range("I1").Value  = "Universal"
$ws.Range("AF1").Value = "Donor_Aid_Budget"
$ws.Range("AG1").Value = "Total_Funds"
$ws.Range("AH1").Value = "%_Private_Funds"
$ws.Range("AI1").Value = "%_MAE_Funds"
$ws.Range("AM1").Value = "Delegation"

# --- Column C numeric updates (rows 2-69) ---
$ws.Range("C2").Value = 2934.187009790061
$ws.Range("C3").Value = 2870.311589353206
$ws.Range("C4").Value = 1873.394108966653
$ws.Range("C5").Value = 4729.735976516416
$ws.Range("C6").Value = 8573.70829393744
$ws.Range("C7").Value = 1286.515571617672
$ws.Range("C8").Value = 2812.435974421079
$ws.Range("C9").Value = 951.6879611168786
$ws.Range("C10").Value = 2983.242707849043
$ws.Range("C11").Value = 2898.942214704482
$ws.Range("C12").Value = 665.6274194933962
$ws.Range("C13").Value = 1904.346464968814
$ws.Range("C14").Value = 5082.354756663512
$ws.Range("C15").Value = 9070.488252857474
$ws.Range("C16").Value = 1303.425880277445
$ws.Range("C17").Value = 2828.483778716848
$ws.Range("C18").Value = 3083.80337578809
$ws.Range("C19").Value = 2965.153206179127
$ws.Range("C20").Value = 1939.33862702996
$ws.Range("C21").Value = 5360.226632400601
$ws.Range("C22").Value = 9603.24465538973
$ws.Range("C23").Value = 3156.723844635973
$ws.Range("C24").Value = 1982.009737844954
$ws.Range("C25").Value = 2999.422762626143
$ws.Range("C26").Value = 5122.180090208862
$ws.Range("C27").Value = 5642.578115155247
$ws.Range("C28").Value = 701.4459636783288
$ws.Range("C29").Value = 11745.7759262897
$ws.Range("C30").Value = 2860.874335573629
$ws.Range("C31").Value = 1325.930225429421
$ws.Range("C32").Value = 1000.829216794104
$ws.Range("C33").Value = 1469.177610078392
$ws.Range("C34").Value = 8629.143035230249
$ws.Range("C35").Value = 3212.740625904757
$ws.Range("C36").Value = 2000.792448761861
$ws.Range("C37").Value = 3056.152683606517
$ws.Range("C38").Value = 5295.682695961288
$ws.Range("C39").Value = 5919.20956823756
$ws.Range("C40").Value = 720.7128711178943
$ws.Range("C41").Value = 11993.48398487312
$ws.Range("C42").Value = 2887.250212489506
$ws.Range("C43").Value = 1360.10887014004
$ws.Range("C44").Value = 1032.277326842402
$ws.Range("C45").Value = 1544.619247249133
$ws.Range("C46").Value = 8965.648525048287
$ws.Range("C47").Value = 11951.20944634967
$ws.Range("C48").Value = 1401.753174264641
$ws.Range("C49").Value = 3008.669179463094
$ws.Range("C50").Value = 1379.14068216006
$ws.Range("C51").Value = 3252.634165082374
$ws.Range("C52").Value = 3137.260298393558
$ws.Range("C53").Value = 2025.814194788851
$ws.Range("C54").Value = 1640.18070024053
$ws.Range("C55").Value = 1060.095015975378
$ws.Range("C56").Value = 11431.15448084494
$ws.Range("C57").Value = 1441.783971398429
$ws.Range("C58").Value = 3012.536723186288
$ws.Range("C59").Value = 1463.71052702022
$ws.Range("C60").Value = 3314.741082534716
$ws.Range("C61").Value = 3210.869677115934
$ws.Range("C62").Value = 2067.29003376698
$ws.Range("C63").Value = 1751.664428859304
$ws.Range("C64").Value = 1093.134170274031
$ws.Range("C65").Value = 1469.192636109792
$ws.Range("C66").Value = 1529.507453727912
$ws.Range("C67").Value = 3382.563653843273
$ws.Range("C68").Value = 3242.636921959078
$ws.Range("C69").Value = 1129.713195979213

# --- Column AL ("Colony") flag updates for the Haiti rows ---
$ws.Range("AL12").Value = 1
$ws.Range("AL28").Value = 1
$ws.Range("AL40").Value = 1
